# Regenerate the "K" (strikeouts) column (column G) in the save_data
# worksheet with freshly computed values, replacing the previous
# Strike# derived figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K value, taken from the regenerated
# s_vals computation (std/mean recalculated upstream).
$kValues = [ordered]@{
    2  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 1
    9  = 2
    10 = 0
    11 = 0
    12 = 1
    13 = 1
    14 = 0
    15 = 0
    16 = 1
    17 = 1
    18 = 5
    19 = 2
    20 = 2
    21 = 0
    22 = 2
    23 = 1
    24 = 1
    25 = 1
    26 = 0
    27 = 0
    28 = 1
    29 = 1
    30 = 0
    31 = 2
    33 = 2
    34 = 3
    35 = 2
    36 = 1
    37 = 6
    38 = 3
    42 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
